# daily auto push: 2026-02-19 07:14 UTC
# Insert one new data row at row 831 (a 2026/02/19 14:00 entry that was
# missing), pushing the existing rows 831:872 down to 832:873.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 831:872 down to 832:873, leaving a blank row 831 behind.
$ws.Rows(831).Insert()

# Fill the newly inserted row 831 with the new record.
# Column A holds dates formatted/stored as plain text ("YYYY/MM/DD"),
# never as real Excel date serials, so force text with a leading
# apostrophe and then strip the resulting quote-prefix style so the
# cell ends up with no explicit style, just like its neighbours.
$ws.Range("A831").Value = "'2026/02/19"
$ws.Range("A831").Style = "Normal"
$ws.Range("B831").Value = "木"
$ws.Range("C831").Value = 14
$ws.Range("D831").Value = 201
